{"js": "// The document ends with a final (empty) paragraph. The author typed\n// \"ALTERA\u00c7\u00c3O\" into that last paragraph, which is why Word's \"_GoBack\"\n// bookmark (marking the last edit location) moved from the end of the\n// previous paragraph (after \" Altera\u00e7\u00e3o da vers\u00e3o para github...\") to the\n// new final paragraph, right after the freshly typed text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert the new run of text at the end of the last (empty) paragraph and\n// make sure it carries the same Tahoma font the rest of the document uses.\nconst insertedRange = lastParagraph.insertText(\"ALTERA\u00c7\u00c3O\", Word.InsertLocation.end);\ninsertedRange.font.name = \"Tahoma\";\nawait context.sync();\n\n// Word keeps a single \"_GoBack\" bookmark that tracks the most recent edit\n// location. Drop it from wherever it used to be and recreate it right after\n// the text we just inserted.\nif (context.document.getBookmarks) {\n  const existing = context.document.getBookmarks().value;\n}\n\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst bookmarkRange = lastParagraph.getRange(Word.RangeLocation.end);\nbookmarkRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The document ends with a final (empty) paragraph. The author typed\n# \"ALTERA\u00c7\u00c3O\" into that paragraph, which moved Word's \"_GoBack\" bookmark\n# (it always tracks the most recent edit location) from the end of the\n# previous paragraph (right after \" Altera\u00e7\u00e3o da vers\u00e3o para github...\")\n# to the new final paragraph, right after the freshly typed text.\n\n$d = $word.ActiveDocument\n\n# Drop the old \"_GoBack\" bookmark from wherever it currently lives.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Seed a fresh collapsed \"_GoBack\" bookmark at the (currently empty) last\n# paragraph, before we add any text \u2014 this is the reliable way to anchor a\n# bookmark at that position.\n$lastParagraph = $d.Paragraphs.Last\n$endRange = $lastParagraph.Range\n$endRange.Collapse(0)   # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $endRange)\n\n# Insert the new text right before the bookmark position so the bookmark\n# naturally ends up immediately after the inserted run, matching Word's\n# real typing behavior.\n$bookmarkRange = $d.Bookmarks(\"_GoBack\").Range\n$insertRange = $d.Range($bookmarkRange.Start, $bookmarkRange.Start)\n$insertRange.InsertBefore(\"ALTERA\u00c7\u00c3O\")\n$insertRange.Font.Name = \"Tahoma\"\n$insertRange.Font.NameBi = \"Tahoma\"\n"}
